# Auto-generated script applying the cryptos.xlsx data refresh described in the commit diff.
# Rows 2-51 hold one crypto-currency entry each (columns B=Coin, C=Link, D=Price, E=Volume(1h)).
# This run updates the Price/Volume figures scraped for this run, and reorders two pairs of rows
# (InjectiveProtocol/Kaspa and PEPE/OKB swapped rank position).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.282.03"
$ws.Range("E2").Value = "  +7.11%  "
$ws.Range("D3").Value = "3.627.71"
$ws.Range("E3").Value = "  +3.91%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'419.65"
$ws.Range("E5").Value = "  +1.18%  "
$ws.Range("D6").Value = "'130.09"
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("D7").Value = "'0.658"
$ws.Range("E7").Value = "  +3.77%  "
$ws.Range("D8").Value = "3.615.39"
$ws.Range("E8").Value = "  +3.78%  "
$ws.Range("E9").Value = "  -0.06%  "
$ws.Range("E10").Value = "  +1.84%  "
$ws.Range("E11").Value = "  +23.98%  "
$ws.Range("D12").Value = "'0.0000429"
$ws.Range("E12").Value = "  +88.84%  "
$ws.Range("D13").Value = "'42.02"
$ws.Range("E13").Value = "  -1.67%  "
$ws.Range("E14").Value = "  +1.46%  "
$ws.Range("D15").Value = "4.204.79"
$ws.Range("E15").Value = "  +4.02%  "
$ws.Range("E16").Value = "  -0.10%  "
$ws.Range("D17").Value = "3.623.47"
$ws.Range("E17").Value = "  +3.68%  "
$ws.Range("D18").Value = "'20.07"
$ws.Range("E18").Value = "  -1.58%  "
$ws.Range("E19").Value = "  +1.84%  "
$ws.Range("D20").Value = "68.189.03"
$ws.Range("E20").Value = "  +7.14%  "
$ws.Range("D21").Value = "'12.39"
$ws.Range("E21").Value = "  -0.54%  "
$ws.Range("D22").Value = "'461.52"
$ws.Range("E22").Value = "  +0.65%  "
$ws.Range("D23").Value = "'89.12"
$ws.Range("E23").Value = "  -1.57%  "
$ws.Range("D24").Value = "'13.44"
$ws.Range("E24").Value = "  +1.06%  "
$ws.Range("E25").Value = "  -6.35%  "
$ws.Range("D26").Value = "'10.20"
$ws.Range("E26").Value = "  -0.86%  "
$ws.Range("D27").Value = "'35.52"
$ws.Range("E27").Value = "  +6.04%  "
$ws.Range("D28").Value = "'3.27"
$ws.Range("E28").Value = "  -1.55%  "
$ws.Range("D29").Value = "'4.98"
$ws.Range("E29").Value = "  +4.54%  "
$ws.Range("E30").Value = "  +4.98%  "
$ws.Range("D31").Value = "'12.22"
$ws.Range("E31").Value = "  -3.08%  "
$ws.Range("E32").Value = "  +4.38%  "
$ws.Range("E33").Value = "  -3.90%  "
$ws.Range("B34").Value = "Kaspa"
$ws.Range("C34").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D34").Value = "'0.158"
$ws.Range("E34").Value = "  -7.96%  "
$ws.Range("B35").Value = "InjectiveProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D35").Value = "'40.31"
$ws.Range("E35").Value = "  +0.47%  "
$ws.Range("D36").Value = "'0.998"
$ws.Range("E36").Value = "  -0.21%  "
$ws.Range("B37").Value = "OKB"
$ws.Range("C37").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D37").Value = "'56.03"
$ws.Range("E37").Value = "  -2.67%  "
$ws.Range("B38").Value = "PEPE"
$ws.Range("C38").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D38").Value = "0.0₃0799"
$ws.Range("E38").Value = "  +19.19%  "
$ws.Range("D39").Value = "'0.0489"
$ws.Range("E39").Value = "  +0.23%  "
$ws.Range("D40").Value = "'0.149"
$ws.Range("E40").Value = "  +9.57%  "
$ws.Range("D41").Value = "'0.998"
$ws.Range("E41").Value = "  -0.18%  "
$ws.Range("D42").Value = "'148.67"
$ws.Range("E42").Value = "  +1.99%  "
$ws.Range("D43").Value = "'2.73"
$ws.Range("E43").Value = "  -3.02%  "
$ws.Range("D44").Value = "'2.93"
$ws.Range("E44").Value = "  -4.46%  "
$ws.Range("E45").Value = "  -3.28%  "
$ws.Range("E46").Value = "  -6.18%  "
$ws.Range("E47").Value = "  +19.77%  "
$ws.Range("E48").Value = "  +9.62%  "
$ws.Range("E49").Value = "  -4.37%  "
$ws.Range("E50").Value = "  -2.99%  "
$ws.Range("D51").Value = "'2.64"
